$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1) values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2: value moved from D2 to C2 with updated value; clear old D2
$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 26.835524999093739

# Row 3: clear B3 and C3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the selection to match the new active range
$ws.Range("B1:E3").Select()
